# Add a new "Consumer complaints" / COM row to the
# iApply_Integration_Overall_Stat sheet (row 41), and update the
# selection on that sheet, per commit:
#   "adding customer complaints form related sql files into repo"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("iApply_Integration_Overall_Stat")
$ws2 = $wb.Worksheets.Item("sp_view_prod_updation")

# 1) Seed row 41 with the same formatting as row 40 (border/style pattern
#    used throughout the table for columns A and D..P).
$ws1.Range("A40:P40").Copy()
$ws1.Range("A41:P41").PasteSpecial(-4122)

# 2) Columns B and C of this new row use slightly different (plain
#    bordered) cell styles than row 40 - pick them up from existing
#    cells elsewhere in the workbook that already use that formatting.
$ws2.Range("D23").Copy()
$ws1.Range("B41").PasteSpecial(-4122)
$ws2.Range("D32").Copy()
$ws1.Range("C41").PasteSpecial(-4122)

# 3) Fill in the row's values. Order is chosen so any newly introduced
#    shared strings land in the same sequence as the source edit.
$ws1.Range("C41").Value = "Consumer complaints"
$ws1.Range("B41").Value = "COM"
$ws1.Range("D41:M41").Value = "Work in Progress"
$ws1.Range("A41").Value = 39
$ws1.Range("N41:P41").Value = "-"

# 4) Update the active selection on the overview sheet.
$ws1.Activate()
$ws1.Range("M38").Select()
